# LOQ4049.xlsx update
# 1) "Semestre ideal" value changes from EQD-9,EQN-12 to EQD-9,EQN-11
# 2) "Requisitos" list drops the LOM3081 and LOQ4054 entries, shifting the
#    remaining two requirements (LOQ4002, LOQ4086) up into rows 24-25 and
#    removing the now-empty rows 26-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "Semestre ideal" (row 9, columns B and C share the same text)
$ws.Range("B9:C9").Value = "EQD-9,EQN-11"

# 2) Rewrite the two requirement rows that survive, then drop the
#    trailing rows that held the removed requirements.
$req1 = "LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n"
$req2 = "LOQ4086 -  Operações Unitárias II  (Requisito fraco)`n"

$ws.Range("B24:C24").Value = $req1
$ws.Range("B25:C25").Value = $req2

$ws.Rows("26:27").Delete()
